$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "strain" values (column F) for data rows 2..43, one per row (3 replicates share a value)
$strains = @(
    "KN99allpha","KN99allpha","KN99allpha",
    "TDY1974","TDY1974","TDY1974",
    "TDY1966","TDY1966","TDY1966",
    "TDY1954","TDY1954","TDY1954",
    "TDY1452","TDY1452","TDY1452",
    "TDY1951","TDY1951","TDY1951",
    "TDY2020","TDY2020","TDY2020",
    "KN99allpha","KN99allpha","KN99allpha",
    "TDY2011","TDY2011","TDY2011",
    "TDY1969","TDY1969","TDY1969",
    "TDY1954","TDY1954","TDY1954",
    "TDY1939","TDY1939","TDY1939",
    "TDY1957","TDY1957","TDY1957",
    "TDY1948","TDY1948","TDY1948"
)

for ($i = 0; $i -lt $strains.Length; $i++) {
    $row = $i + 2

    # Column B (harvester): "Retrofitted_1874" -> "S.GISH" for every data row
    $ws.Cells.Item($row, 2).Value = "S.GISH"

    # Column D (experimentDesign): new column populated with "90minuteInduction"
    $ws.Cells.Item($row, 4).Value = "90minuteInduction"

    # Column F (strain): new column populated with the strain name
    $ws.Cells.Item($row, 6).Value = $strains[$i]
}

# Match the author's final selection in the sheet view
$ws.Range("F42:F43").Select() | Out-Null
